$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192847013473511
$ws.Range("B1").Value = 2.240187168121338
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.280818462371826
$ws.Range("E1").Value = 1.214947700500488
